# Import fixed to stop from saving stages that are already in DB
#
# The classifier-stage column headers in row 5 (F5:K5) were pointing at
# stage codes that had already been imported/saved, so the importer now
# skips those and the headers shift forward to the next unseen stages.
# A couple of dependent labels (and a couple of numbers that were
# recomputed once the duplicate-stage rows dropped out of the import)
# move together with that fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# --- Row 5: classifier stage headers shift to the next un-saved stages ---
$ws.Range("F5").Value = "CLC-09"
$ws.Range("G5").Value = "CLC-11"
$ws.Range("H5").Value = "CLC-13"
$ws.Range("I5").Value = "CLC-15"
$ws.Range("J5").Value = "CLC-17"
$ws.Range("K5").Value = "CLC-19"

# --- Note about which stage got dropped for Jarno Virta now references CLC-07 ---
$ws.Range("A16").Value = "Jarno Virta (CLC-07 result dropped out because more than 8 results)"

# --- New result recorded for the newly-imported stage (column K) ---
$ws.Range("K6").Value = 1
$ws.Range("K6").NumberFormat = "0.00"

# --- Percentage table: new value for the same column, row 25 ---
$ws.Range("K25").Value = 0.55556000000000005

# --- Row 26 (Jerry Miculek) percentages: recomputed values used for the
#     TOP 4 average are now highlighted in bold, and the TOP4 avg itself
#     (L26) is updated ---
$ws.Range("B26").Font.Bold = $true
$ws.Range("E26").Font.Bold = $true
$ws.Range("H26").Font.Bold = $true
$ws.Range("K26").Font.Bold = $true
$ws.Range("L26").Value = 0.87544999999999995

# --- Ranking table: Jerry Miculek's ranking percentage recalculated ---
$ws.Range("B33").Value = 0.75329999999999997

# --- Restore the selection left on the sheet after the edit ---
$ws.Range("D14").Select()
